$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New name entries (columns N/O), appended after row 34
$ws.Range("N35").Value = "Baeddan"
$ws.Range("O35").Value = "Баэддан"
$ws.Range("N36").Value = "Bedwyr"
$ws.Range("O36").Value = "Бедвир"
$ws.Range("N37").Value = "Esni"
$ws.Range("O37").Value = "Эсни"
$ws.Range("N38").Value = "Cei"
$ws.Range("O38").Value = "Цей"

# New treasure entries (columns K/L), appended after row 27
$ws.Range("K28").Value = "Stolen Star"
$ws.Range("L28").Value = "Украденная звезда"
$ws.Range("K29").Value = "Javelin of Thorns"
$ws.Range("L29").Value = "Дрот шипов"

$ws.Range("F32").Select()
